$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14 through 28 first (old extra data no longer present in final sheet)
$ws.Range("A14:F28").EntireRow.Delete()

# New data for columns B and C (rows 2-13)
$colB = @(
  "NSE:AEGISCHEM",
  "NSE:AGRITECH",
  "NSE:BHEL",
  "NSE:CONCORDBIO",
  "NSE:IGL",
  "NSE:JTEKTINDIA",
  "NSE:LOKESHMACH",
  "NSE:LTTS",
  "NSE:MAHLIFE",
  "NSE:MAKEINDIA",
  "NSE:MSPL",
  "NSE:PRICOLLTD"
)

$colC = @(
  "NSE:APOLLO",
  "NSE:ASALCBR",
  "NSE:DCMSHRIRAM",
  "NSE:GRAPHITE",
  "NSE:HERCULES",
  "NSE:IDEAFORGE",
  "NSE:JAIPURKURT",
  "NSE:KANPRPLA",
  "NSE:KELLTONTEC",
  "NSE:MRO-TEK",
  "",
  ""
)

for ($i = 0; $i -lt $colB.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $colB[$i]
  $ws.Cells.Item($row, 3).Value = $colC[$i]
}

# Column D: clear all existing values (D2 had NSE:PIIND, rest already empty)
$ws.Range("D2:D13").Value = ""

# Column F: update rows 2-3, clear row 4 (rest already empty)
$ws.Cells.Item(2, 6).Value = "NSE:BHEL"
$ws.Cells.Item(3, 6).Value = "NSE:LTTS"
$ws.Cells.Item(4, 6).Value = ""
